# Add the new daily data (1 June 2021 - 30 June 2021, serial dates
# 44348-44377) as rows 450-479 on every sheet of the workbook, mirroring
# how the source spreadsheet was extended by the original author.
#
# Only row 450 receives data in columns C/D(/E); rows 451-479 only get the
# date in column A (the remaining days had not been reported yet at the
# time of this commit).

$wb = $excel.ActiveWorkbook

# date serials for rows 450 (2021-06-01) .. 479 (2021-06-30)
$firstDate = 44348
$firstRow  = 450
$lastRow   = 479

function Fill-DatesAndFormat($ws) {
    # Copy the date cell's number format down first so every new row
    # matches the existing dd/mm/yyyy styling (style index 5).
    $ws.Range("A449").Copy()
    $ws.Range("A450:A479").PasteSpecial(-4122)

    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $firstDate + ($r - $firstRow)
    }
}

# ---------------------------------------------------------------------
# Sheet 1: "Nuovi casi"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Nuovi casi")
Fill-DatesAndFormat $ws1
$ws1.Range("C449:D449").Copy()
$ws1.Range("C450:D450").PasteSpecial(-4122)
$ws1.Range("C450").Value = 28
$ws1.Range("D450").Formula = "=AVERAGE(C444:C450)"

# ---------------------------------------------------------------------
# Sheet 2: "Deceduti"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Deceduti")
Fill-DatesAndFormat $ws2
$ws2.Range("C449:D449").Copy()
$ws2.Range("C450:D450").PasteSpecial(-4122)
$ws2.Range("C450").Value = 3
$ws2.Range("D450").Formula = "=AVERAGE(C444:C450)"

# ---------------------------------------------------------------------
# Sheet 3: "Dimessi   Guariti"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Dimessi   Guariti")
Fill-DatesAndFormat $ws3
$ws3.Range("C449:D449").Copy()
$ws3.Range("C450:D450").PasteSpecial(-4122)
$ws3.Range("C450").Value = 126
$ws3.Range("D450").Formula = "=AVERAGE(C444:C450)"

# ---------------------------------------------------------------------
# Sheet 4: "Ricoveri"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Ricoveri")
Fill-DatesAndFormat $ws4
$ws4.Range("C449:E449").Copy()
$ws4.Range("C450:E450").PasteSpecial(-4122)
$ws4.Range("C450").Value = 93
$ws4.Range("D450").Formula = "=AVERAGE(C444:C450)"
$ws4.Range("E450").Formula = "=C450-C449"

# ---------------------------------------------------------------------
# Sheet 5: "Terapia"
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Terapia")
Fill-DatesAndFormat $ws5
$ws5.Range("C449:D449").Copy()
$ws5.Range("C450:D450").PasteSpecial(-4122)
$ws5.Range("C450").Value = 6
$ws5.Range("D450").Formula = "=AVERAGE(C444:C450)"

# ---------------------------------------------------------------------
# Restore the view state: "Ricoveri" (sheet index 4) was the active tab,
# each sheet's selection now sits on the newly appended data.
# ---------------------------------------------------------------------
$ws1.Range("A450:A479").Select()
$ws2.Range("A450:D450").Select()
$ws3.Range("A450:D450").Select()
$ws5.Range("A450").Select()
$ws4.Activate()
$ws4.Range("A450:D450").Select()

Write-Output "Added rows 450-479 to all sheets"
